# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.107.83"
$ws.Range("E2").Value = "'  +0.41%  "
$ws.Range("D3").Value = "'3.618.83"
$ws.Range("E3").Value = "'  +3.56%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'604.87"
$ws.Range("E5").Value = "'  +0.43%  "
$ws.Range("D6").Value = "'195.68"
$ws.Range("E6").Value = "'  -0.99%  "
$ws.Range("E7").Value = "'  +0.46%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E9").Value = "'  -1.35%  "
$ws.Range("D10").Value = "'0.651"
$ws.Range("E10").Value = "'  -0.13%  "
$ws.Range("D11").Value = "'53.93"
$ws.Range("E11").Value = "'  -0.30%  "
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("E12").Value = "'  +0.92%  "
$ws.Range("D13").Value = "'9.54"
$ws.Range("E13").Value = "'  +0.02%  "
$ws.Range("D14").Value = "'4.189.96"
$ws.Range("E14").Value = "'  +3.33%  "
$ws.Range("D15").Value = "'13.27"
$ws.Range("E15").Value = "'  +5.25%  "
$ws.Range("D16").Value = "'592.33"
$ws.Range("E16").Value = "'  -0.08%  "
$ws.Range("E17").Value = "'  +1.70%  "
$ws.Range("D18").Value = "'70.298.10"
$ws.Range("E18").Value = "'  +0.57%  "
$ws.Range("D19").Value = "'3.593.15"
$ws.Range("E19").Value = "'  +2.72%  "
$ws.Range("E20").Value = "'  +1.71%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "'  +1.35%  "
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'5.16"
$ws.Range("E23").Value = "'  +3.36%  "
$ws.Range("D24").Value = "'102.94"
$ws.Range("E24").Value = "'  -0.69%  "
$ws.Range("E25").Value = "'  +0.99%  "
$ws.Range("D26").Value = "'3.07"
$ws.Range("E26").Value = "'  -0.34%  "
$ws.Range("D27").Value = "'10.88"
$ws.Range("E27").Value = "'  -0.46%  "
$ws.Range("D28").Value = "'9.61"
$ws.Range("E28").Value = "'  -1.49%  "
$ws.Range("D29").Value = "'33.96"
$ws.Range("B30").Value = "'dogwifhat"
$ws.Range("C30").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D30").Value = "'4.40"
$ws.Range("E30").Value = "'  -2.08%  "
$ws.Range("B31").Value = "'NEARProtocol"
$ws.Range("C31").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.13"
$ws.Range("E31").Value = "'  -1.73%  "
$ws.Range("D32").Value = "'12.34"
$ws.Range("E32").Value = "'  -2.76%  "
$ws.Range("E33").Value = "'  +0.91%  "
$ws.Range("D34").Value = "'63.24"
$ws.Range("E34").Value = "'  -0.47%  "
$ws.Range("D35").Value = "'3.931.35"
$ws.Range("E35").Value = "'  +6.72%  "
$ws.Range("D36").Value = "'0.0₃0861"
$ws.Range("E36").Value = "'  +7.45%  "
$ws.Range("D37").Value = "'3.19"
$ws.Range("E37").Value = "'  +7.13%  "
$ws.Range("D38").Value = "'524.61"
$ws.Range("E38").Value = "'  +2.06%  "
$ws.Range("E39").Value = "'  +0.19%  "
$ws.Range("D40").Value = "'37.16"
$ws.Range("E40").Value = "'  +1.78%  "
$ws.Range("E41").Value = "'  +0.97%  "
$ws.Range("D42").Value = "'3.54"
$ws.Range("E42").Value = "'  +0.63%  "
$ws.Range("E43").Value = "'  -2.28%  "
$ws.Range("E44").Value = "'  +0.15%  "
$ws.Range("B45").Value = "'ThetaToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "'  +0.88%  "
$ws.Range("B46").Value = "'Stellar"
$ws.Range("C46").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.141"
$ws.Range("E46").Value = "'  +1.15%  "
$ws.Range("B47").Value = "'ApeXProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "'  +0.12%  "
$ws.Range("D48").Value = "'8.62"
$ws.Range("E48").Value = "'  -1.38%  "
$ws.Range("E49").Value = "'  -0.09%  "
$ws.Range("D50").Value = "'0.000247"
$ws.Range("E50").Value = "'  +3.00%  "
$ws.Range("D51").Value = "'1.31"
$ws.Range("E51").Value = "'  +3.19%  "
